$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3.476304544167619
$ws.Range("D2").Value = 3.342861255556732
$ws.Range("E2").Value = 40.47176262065683
$ws.Range("F2").Value = 16.7669461497492
$ws.Range("G2").Value = 3.571768659760583
$ws.Range("I2").Value = 14.39738836212795
$ws.Range("M2").Value = 57.91617868534593
$ws.Range("O2").Value = 14.76877119249959
$ws.Range("C3").Value = 3.334953336422328
$ws.Range("D3").Value = 3.221367391710547
$ws.Range("E3").Value = 37.75138303523279
$ws.Range("F3").Value = 16.95767206342498
$ws.Range("G3").Value = 3.575884284977013
$ws.Range("I3").Value = 14.60409137007161
$ws.Range("M3").Value = 54.35550820080955
$ws.Range("O3").Value = 15.04093699823042
$ws.Range("C4").Value = 3.244691757843198
$ws.Range("D4").Value = 3.171059348707021
$ws.Range("E4").Value = 35.98041707787855
$ws.Range("F4").Value = 17.09384532706991
$ws.Range("G4").Value = 3.578509361899191
$ws.Range("I4").Value = 14.74381313235523
$ws.Range("M4").Value = 52.04303753748426
$ws.Range("O4").Value = 15.22196686562206
$ws.Range("C5").Value = 3.207074058257863
$ws.Range("D5").Value = 3.167217141609902
$ws.Range("E5").Value = 35.23326466210144
$ws.Range("F5").Value = 17.15395081477577
$ws.Range("G5").Value = 3.579603990821591
$ws.Range("I5").Value = 14.80386814740245
$ws.Range("M5").Value = 51.06898535093034
$ws.Range("O5").Value = 15.29912243925746
$ws.Range("C6").Value = 3.200778417950994
$ws.Range("D6").Value = 3.166594312008303
$ws.Range("E6").Value = 35.10765763722933
$ws.Range("F6").Value = 17.16420472170481
$ws.Range("G6").Value = 3.579787262867508
$ws.Range("I6").Value = 14.81402559352947
$ws.Range("M6").Value = 50.90533207952289
$ws.Range("O6").Value = 15.31213523491345
$ws.Range("C7").Value = 3.244187761188732
$ws.Range("D7").Value = 3.171006515235379
$ws.Range("E7").Value = 35.97044405549916
$ws.Range("F7").Value = 17.09463748618448
$ws.Range("G7").Value = 3.578524023376547
$ws.Range("I7").Value = 14.74461056604613
$ws.Range("M7").Value = 52.03002940477588
$ws.Range("O7").Value = 15.22299385656659
$ws.Range("C8").Value = 3.428307985960186
$ws.Range("D8").Value = 3.301583366345975
$ws.Range("E8").Value = 39.55456471346349
$ws.Range("F8").Value = 16.82864160109181
$ws.Range("G8").Value = 3.573167515019288
$ws.Range("I8").Value = 14.46593865959478
$ws.Range("M8").Value = 56.71465714749222
$ws.Range("O8").Value = 14.85965323867943
$ws.Range("C9").Value = 3.760439405120442
$ws.Range("D9").Value = 3.587951167072841
$ws.Range("E9").Value = 45.79201205668925
$ws.Range("F9").Value = 16.46649431085119
$ws.Range("G9").Value = 3.563430537808864
$ws.Range("I9").Value = 14.02589357867655
$ws.Range("M9").Value = 64.90092857163805
$ws.Range("O9").Value = 14.26315838900036
$ws.Range("C10").Value = 3.985299899419243
$ws.Range("D10").Value = 3.783070663488931
$ws.Range("E10").Value = 49.90317859093299
$ws.Range("F10").Value = 16.30887888768651
$ws.Range("G10").Value = 3.556728203287531
$ws.Range("I10").Value = 13.77442895059257
$ws.Range("M10").Value = 70.31005263774509
$ws.Range("O10").Value = 13.90392484520304
$ws.Range("C11").Value = 4.08319656546668
$ws.Range("D11").Value = 3.868405638425283
$ws.Range("E11").Value = 51.67346889593462
$ws.Range("F11").Value = 16.26311145027135
$ws.Range("G11").Value = 3.553773390399056
$ws.Range("I11").Value = 13.67722456640414
$ws.Range("M11").Value = 72.64094891128616
$ws.Range("O11").Value = 13.7596413062023
$ws.Range("C12").Value = 4.119620374759179
$ws.Range("D12").Value = 3.900220695137537
$ws.Range("E12").Value = 52.32966520989347
$ws.Range("F12").Value = 16.24970087712955
$ws.Range("G12").Value = 3.5526676925367
$ws.Range("I12").Value = 13.64302550826068
$ws.Range("M12").Value = 73.50509775027406
$ws.Range("O12").Value = 13.70793473025431
$ws.Range("C13").Value = 4.111804873721499
$ws.Range("D13").Value = 3.893391052122997
$ws.Range("E13").Value = 52.18896889314987
$ws.Range("F13").Value = 16.252411732261
$ws.Range("G13").Value = 3.552905241216475
$ws.Range("I13").Value = 13.65027254989014
$ws.Range("M13").Value = 73.31980858910056
$ws.Range("O13").Value = 13.71893736715893
$ws.Range("C14").Value = 4.08620623119143
$ws.Range("D14").Value = 3.871033123324953
$ws.Range("E14").Value = 51.72773751985636
$ws.Range("F14").Value = 16.26192850511949
$ws.Range("G14").Value = 3.553682160730086
$ws.Range("I14").Value = 13.67435786967575
$ws.Range("M14").Value = 72.71241303992181
$ws.Range("O14").Value = 13.75532758009622
$ws.Range("C15").Value = 4.070441579865666
$ws.Range("D15").Value = 3.857273038284782
$ws.Range("E15").Value = 51.44337985050839
$ws.Range("F15").Value = 16.26827395773
$ws.Range("G15").Value = 3.554159759323384
$ws.Range("I15").Value = 13.68945492451533
$ws.Range("M15").Value = 72.33795971066056
$ws.Range("O15").Value = 13.77800476849579
$ws.Range("C16").Value = 3.978812296516062
$ws.Range("D16").Value = 3.777424253193506
$ws.Range("E16").Value = 49.78549194973996
$ws.Range("F16").Value = 16.31240822538049
$ws.Range("G16").Value = 3.556923175121765
$ws.Range("I16").Value = 13.78113909350044
$ws.Range("M16").Value = 70.1551243486163
$ws.Range("O16").Value = 13.91375419320501
$ws.Range("C17").Value = 3.921463218154494
$ws.Range("D17").Value = 3.727556622303853
$ws.Range("E17").Value = 48.74297915822704
$ws.Range("F17").Value = 16.34626308778967
$ws.Range("G17").Value = 3.558642331957822
$ws.Range("I17").Value = 13.84188378111359
$ws.Range("M17").Value = 68.78289284747197
$ws.Range("O17").Value = 14.00205664307449
$ws.Range("C18").Value = 3.888064740245451
$ws.Range("D18").Value = 3.69855188417306
$ws.Range("E18").Value = 48.13393857097529
$ws.Range("F18").Value = 16.36816523261235
$ws.Range("G18").Value = 3.559640023207515
$ws.Range("I18").Value = 13.87843176859088
$ws.Range("M18").Value = 67.98139961853012
$ws.Range("O18").Value = 14.05463543290603
$ws.Range("C19").Value = 3.87668622070631
$ws.Range("D19").Value = 3.688676284065928
$ws.Range("E19").Value = 47.92610514571452
$ws.Range("F19").Value = 16.37599230517209
$ws.Range("G19").Value = 3.559979358184159
$ws.Range("I19").Value = 13.89107853214795
$ws.Range("M19").Value = 67.70792454204363
$ws.Range("O19").Value = 14.07273978558993
$ws.Range("C20").Value = 3.927610973026165
$ws.Range("D20").Value = 3.732898518989707
$ws.Range("E20").Value = 48.85492945115784
$ws.Range("F20").Value = 16.34240630752585
$ws.Range("G20").Value = 3.558458408145516
$ws.Range("I20").Value = 13.83524986496313
$ws.Range("M20").Value = 68.93023316801961
$ws.Range("O20").Value = 13.99247022138893
$ws.Range("C21").Value = 4.093742857435595
$ws.Range("D21").Value = 3.877613787449467
$ws.Range("E21").Value = 51.86359548839178
$ws.Range("F21").Value = 16.25902529924958
$ws.Range("G21").Value = 3.553453604229646
$ws.Range("I21").Value = 13.66721146396436
$ws.Range("M21").Value = 72.89132096063811
$ws.Range("O21").Value = 13.74455787982243
$ws.Range("C22").Value = 4.198540085115629
$ws.Range("D22").Value = 3.969280747956073
$ws.Range("E22").Value = 53.74740705571445
$ws.Range("F22").Value = 16.22746711066646
$ws.Range("G22").Value = 3.550259609969122
$ws.Range("I22").Value = 13.57266970209716
$ws.Range("M22").Value = 75.37230562142062
$ws.Range("O22").Value = 13.59970452162914
$ws.Range("C23").Value = 4.142958100137584
$ws.Range("D23").Value = 3.920624593256936
$ws.Range("E23").Value = 52.74946791505121
$ws.Range("F23").Value = 16.24215006034431
$ws.Range("G23").Value = 3.551957371315405
$ws.Range("I23").Value = 13.621683353439
$ws.Range("M23").Value = 74.05796939377237
$ws.Range("O23").Value = 13.67538233196495
$ws.Range("C24").Value = 3.924832904101109
$ws.Range("D24").Value = 3.730484490945673
$ws.Range("E24").Value = 48.80434690817094
$ws.Range("F24").Value = 16.34414237699092
$ws.Range("G24").Value = 3.558541531032786
$ws.Range("I24").Value = 13.83824400720335
$ws.Range("M24").Value = 68.86365980087406
$ws.Range("O24").Value = 13.99679861041997
$ws.Range("C25").Value = 3.673859722901696
$ws.Range("D25").Value = 3.513116742077064
$ws.Range("E25").Value = 44.18820805263049
$ws.Range("F25").Value = 16.54623268412135
$ws.Range("G25").Value = 3.565984152374167
$ws.Range("I25").Value = 14.1328446843305
$ws.Range("M25").Value = 62.79307609995573
$ws.Range("O25").Value = 14.76315838900036
